$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the two new sheets (gemini-pro, gemini-flash) at the end, in
#    that order, matching the sheetId / rId ordering in the target.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsGeminiPro = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsGeminiPro.Name = "gemini-pro"

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsGeminiFlash = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)
$wsGeminiFlash.Name = "gemini-flash"

$rounds = @("Round 0","Round 1","Round 2","Round 3","Round 4","Round 5","Round 6","Round 7","Round 8","Round 9","Round 10")

# ---------------------------------------------------------------------
# 2. Populate gemini-flash FIRST so the "Everyone bids 2" shared string
#    lands before the avg/median/std.dev labels (matches shared-string
#    table ordering in the target: 14=Everyone bids 2, 15=avg,
#    16=median, 17=std.dev).
# ---------------------------------------------------------------------
for ($i = 0; $i -lt $rounds.Length; $i++) {
    $r = $i + 2
    $wsGeminiFlash.Cells.Item($r, 1).Value = $rounds[$i]
    $wsGeminiFlash.Cells.Item($r, 2).Value = "Everyone bids 2"
}

# ---------------------------------------------------------------------
# 3. Populate gemini-pro data rows (Round label / bid 8 / paid 38) then
#    the avg / median / std.dev summary rows.
# ---------------------------------------------------------------------
for ($i = 0; $i -lt $rounds.Length; $i++) {
    $r = $i + 2
    $wsGeminiPro.Cells.Item($r, 1).Value = $rounds[$i]
    $wsGeminiPro.Cells.Item($r, 2).Value = 8
    $wsGeminiPro.Cells.Item($r, 3).Value = 38
}

$wsGeminiPro.Range("C13").Value = "avg"
$wsGeminiPro.Range("D13").Formula = "=AVERAGE(C1:C12)"
$wsGeminiPro.Range("C14").Value = "median"
$wsGeminiPro.Range("D14").Formula = "=MEDIAN(C1:C12)"
$wsGeminiPro.Range("C15").Value = "std.dev"
$wsGeminiPro.Range("D15").Formula = "=STDEV.S(C1:C12)"

$wsGeminiPro.Range("D19").Select() | Out-Null
$wsGeminiFlash.Range("B16").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. Add the avg / median / std.dev summary formulas to the three
#    existing sheets.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("D14").Formula = "=AVERAGE(C2:C13)"
$ws1.Range("D15").Formula = "=MEDIAN(C2:C13)"
$ws1.Range("D16").Formula = "=STDEV.S(C2:C13)"
$ws1.Range("D14:D16").Select() | Out-Null

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("D13").Formula = "=AVERAGE(C1:C12)"
$ws2.Range("D14").Formula = "=MEDIAN(C1:C12)"
$ws2.Range("D15").Formula = "=STDEV.S(C1:C12)"
$ws2.Range("D13:D15").Select() | Out-Null

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("D13").Formula = "=AVERAGE(C1:C12)"
$ws3.Range("D14").Formula = "=MEDIAN(C1:C12)"
$ws3.Range("D15").Formula = "=STDEV.S(C1:C12)"
$ws3.Range("D21").Select() | Out-Null

# ---------------------------------------------------------------------
# 5. Activate "gemini-pro" as the selected tab, matching activeTab=3
#    (0-based index of the 4th sheet) and tabSelected moving there.
# ---------------------------------------------------------------------
$wsGeminiPro.Activate() | Out-Null
